{"js": "// Remove the \"ElencoCommenti\" glossary row from the terminology table\n// (\"Nome / Tipo / Descrizione\"). This matches the commit\n// \"Rimosso Elenco Commenti\" \u2014 the row describing the ElencoCommenti\n// entity is deleted while all other rows stay intact.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Find the terminology table: its first row starts with \"Nome\" and it\n// also contains a row whose first cell is exactly \"ElencoCommenti\".\nlet targetRow = null;\n\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    if (row.cells.items.length > 0) {\n      row.cells.items[0].body.load(\"text\");\n    }\n  }\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    if (row.cells.items.length > 0 && row.cells.items[0].body.text.trim() === \"ElencoCommenti\") {\n      targetRow = row;\n      break;\n    }\n  }\n  if (targetRow) break;\n}\n\nif (targetRow) {\n  targetRow.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the \"ElencoCommenti\" glossary row from the terminology table\n# (\"Nome / Tipo / Descrizione\"). This matches the commit\n# \"Rimosso Elenco Commenti\" - the row describing the ElencoCommenti\n# entity is deleted while all other rows stay intact.\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    $n = $t.Rows.Count\n    for ($i = $n; $i -ge 1; $i--) {\n        $row = $t.Rows.Item($i)\n        # Cell range text includes the trailing cell-mark (CR + cell marker),\n        # strip those control characters before comparing.\n        $cellText = $row.Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7).Trim()\n        if ($cellText -eq \"ElencoCommenti\") {\n            $row.Delete()\n        }\n    }\n}\n"}
